$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 17545988
$ws.Range("I100").Value = 22223348
$ws.Range("K100").Value = 22223348
$ws.Range("M100").Value = -22222807

$ws.Range("H113").Value = 2438.4827
$ws.Range("I113").Value = 1788.4375
$ws.Range("J113").Value = 3238.5386
$ws.Range("K113").Value = 1788.4375
$ws.Range("L113").Value = 3238.5386
$ws.Range("M113").Value = 1465.5625
$ws.Range("N113").Value = -9746.5386

$ws.Range("H137").Value = 1094.36
$ws.Range("I137").Value = 813.6842
$ws.Range("J137").Value = 1983.1666
$ws.Range("K137").Value = 2441.0526
$ws.Range("L137").Value = 5949.4998
$ws.Range("M137").Value = 108.9474
$ws.Range("N137").Value = -11049.4998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1781.2963
$ws.Range("I2").Value = 916.8
$ws.Range("J2").Value = 2861.9167
$ws.Range("K2").Value = 916.8
$ws.Range("L2").Value = 2861.9167
$ws.Range("M2").Value = -803.8
$ws.Range("N2").Value = -3087.9167

$ws.Range("H32").Value = 450784.72
$ws.Range("I32").Value = 4303.129
$ws.Range("J32").Value = 14291714
$ws.Range("K32").Value = 4303.129
$ws.Range("L32").Value = 14291714
$ws.Range("M32").Value = -4016.129
$ws.Range("N32").Value = -14292288

$ws.Range("H45").Value = 3646.8635
$ws.Range("I45").Value = 3913.0908
$ws.Range("J45").Value = 3380.6365
$ws.Range("K45").Value = 3913.0908
$ws.Range("L45").Value = 3380.6365
$ws.Range("M45").Value = -3536.0908
$ws.Range("N45").Value = -4134.636500000001

$ws.Range("H61").Value = 2122.742
$ws.Range("I61").Value = 2061.5417
$ws.Range("J61").Value = 2332.5715
$ws.Range("K61").Value = 2061.5417
$ws.Range("L61").Value = 2332.5715
$ws.Range("M61").Value = -1849.5417
$ws.Range("N61").Value = -2756.5715

$ws.Range("H74").Value = 3227.2144
$ws.Range("I74").Value = 3523.6
$ws.Range("J74").Value = 2486.25
$ws.Range("K74").Value = 3523.6
$ws.Range("L74").Value = 2486.25
$ws.Range("M74").Value = -2649.6
$ws.Range("N74").Value = -4234.25

$ws.Range("H77").Value = 3227.2144
$ws.Range("I77").Value = 3523.6
$ws.Range("J77").Value = 2486.25
$ws.Range("K77").Value = 17618
$ws.Range("L77").Value = 12431.25
$ws.Range("M77").Value = -13250
$ws.Range("N77").Value = -21167.25

$ws.Range("H102").Value = 4317.773
$ws.Range("I102").Value = 3792.5
$ws.Range("J102").Value = 5718.5
$ws.Range("K102").Value = 3792.5
$ws.Range("L102").Value = 5718.5
$ws.Range("M102").Value = -2170.5
$ws.Range("N102").Value = -8962.5

$ws.Range("H110").Value = 1395.9
$ws.Range("I110").Value = 1000.375
$ws.Range("J110").Value = 2978
$ws.Range("K110").Value = 1000.375
$ws.Range("L110").Value = 2978
$ws.Range("M110").Value = 1044.625
$ws.Range("N110").Value = -7068

$ws.Range("H116").Value = 1781.2963
$ws.Range("I116").Value = 916.8
$ws.Range("J116").Value = 2861.9167
$ws.Range("K116").Value = 916.8
$ws.Range("L116").Value = 2861.9167
$ws.Range("M116").Value = 1377.2
$ws.Range("N116").Value = -7449.9167

$ws.Range("H122").Value = 2316.6287
$ws.Range("I122").Value = 2267.8333
$ws.Range("J122").Value = 2609.4
$ws.Range("K122").Value = 6803.499899999999
$ws.Range("L122").Value = 7828.200000000001
$ws.Range("M122").Value = -4353.499899999999
$ws.Range("N122").Value = -12728.2

$ws.Range("H132").Value = 2649.6086
$ws.Range("I132").Value = 2195
$ws.Range("J132").Value = 3937.6667
$ws.Range("K132").Value = 6585
$ws.Range("L132").Value = 11813.0001
$ws.Range("M132").Value = -4055
$ws.Range("N132").Value = -16873.0001

$ws.Range("H136").Value = 2122.742
$ws.Range("I136").Value = 2061.5417
$ws.Range("J136").Value = 2332.5715
$ws.Range("K136").Value = 6184.625100000001
$ws.Range("L136").Value = 6997.7145
$ws.Range("M136").Value = -3634.625100000001
$ws.Range("N136").Value = -12097.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1781.2963
$ws.Range("I3").Value = 916.8
$ws.Range("J3").Value = 2861.9167
$ws.Range("K3").Value = 916.8
$ws.Range("L3").Value = 2861.9167
$ws.Range("M3").Value = -802.8
$ws.Range("N3").Value = -3089.9167

$ws.Range("H44").Value = 11984
$ws.Range("J44").Value = 11984
$ws.Range("L44").Value = 11984

$ws.Range("H99").Value = 1837.9286
$ws.Range("I99").Value = 1001.2857
$ws.Range("J99").Value = 2674.5715
$ws.Range("K99").Value = 1001.2857
$ws.Range("L99").Value = 2674.5715
$ws.Range("M99").Value = 496.7143
$ws.Range("N99").Value = -5670.5715

$ws.Range("H104").Value = 55684
$ws.Range("J104").Value = 55684
$ws.Range("L104").Value = 55684

$ws.Range("H105").Value = 1692.5555
$ws.Range("I105").Value = 1649.8334
$ws.Range("J105").Value = 1704.762
$ws.Range("K105").Value = 1649.8334
$ws.Range("L105").Value = 1704.762
$ws.Range("M105").Value = 97.16660000000002
$ws.Range("N105").Value = -5198.762

$ws.Range("H107").Value = 910568.6
$ws.Range("I107").Value = 1315.1666
$ws.Range("J107").Value = 2001672.8
$ws.Range("K107").Value = 1315.1666
$ws.Range("L107").Value = 2001672.8
$ws.Range("M107").Value = 604.8334
$ws.Range("N107").Value = -2005512.8

$ws.Range("H134").Value = 4636.175
$ws.Range("I134").Value = 1017.4865
$ws.Range("J134").Value = 49266.668
$ws.Range("K134").Value = 3052.4595
$ws.Range("L134").Value = 147800.004
$ws.Range("M134").Value = -517.4594999999999
$ws.Range("N134").Value = -152870.004

$ws.Range("H140").Value = 88133.336
$ws.Range("J140").Value = 88133.336
$ws.Range("L140").Value = 88133.336
$ws.Range("N140").Value = -98493.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 759.4
$ws.Range("J16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("N16").Value = -1574

$ws.Range("H31").Value = 2527.3257
$ws.Range("I31").Value = 1738.6562
$ws.Range("J31").Value = 4821.636
$ws.Range("K31").Value = 1738.6562
$ws.Range("L31").Value = 4821.636
$ws.Range("M31").Value = -1443.6562
$ws.Range("N31").Value = -5411.636

$ws.Range("H34").Value = 2527.3257
$ws.Range("I34").Value = 1738.6562
$ws.Range("J34").Value = 4821.636
$ws.Range("K34").Value = 1738.6562
$ws.Range("L34").Value = 4821.636
$ws.Range("M34").Value = -1536.6562
$ws.Range("N34").Value = -5225.636

$ws.Range("H99").Value = 2642
$ws.Range("I99").Value = 2858.5715
$ws.Range("J99").Value = 2389.3333
$ws.Range("K99").Value = 2858.5715
$ws.Range("L99").Value = 2389.3333
$ws.Range("M99").Value = -1360.5715
$ws.Range("N99").Value = -5385.3333

$ws.Range("H105").Value = 1322
$ws.Range("I105").Value = 1003.3333
$ws.Range("J105").Value = 1800
$ws.Range("K105").Value = 1003.3333
$ws.Range("L105").Value = 1800
$ws.Range("M105").Value = 743.6667
$ws.Range("N105").Value = -5294

$ws.Range("H107").Value = 450.55
$ws.Range("I107").Value = 466.75
$ws.Range("K107").Value = 466.75
$ws.Range("M107").Value = 1453.25

$ws.Range("H113").Value = 759.4
$ws.Range("J113").Value = 1000
$ws.Range("L113").Value = 1000
$ws.Range("N113").Value = -5340

$ws.Range("H122").Value = 1618
$ws.Range("I122").Value = 1397.375
$ws.Range("J122").Value = 2122.2856
$ws.Range("K122").Value = 4192.125
$ws.Range("L122").Value = 6366.8568
$ws.Range("M122").Value = -1742.125
$ws.Range("N122").Value = -11266.8568

$ws.Range("H126").Value = 2642
$ws.Range("I126").Value = 2858.5715
$ws.Range("J126").Value = 2389.3333
$ws.Range("K126").Value = 8575.7145
$ws.Range("L126").Value = 7167.999899999999
$ws.Range("M126").Value = -6105.7145
$ws.Range("N126").Value = -12107.9999

$ws.Range("H134").Value = 1417.0834
$ws.Range("I134").Value = 1409.5454
$ws.Range("K134").Value = 4228.6362
$ws.Range("M134").Value = -1693.6362

$ws.Range("H138").Value = 32500
$ws.Range("J138").Value = 32500
$ws.Range("L138").Value = 32500
$ws.Range("N138").Value = -42780

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2950
$ws.Range("J39").Value = 3833.3333
$ws.Range("L39").Value = 11499.9999
$ws.Range("N39").Value = -12087.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2343.2307
$ws.Range("I102").Value = 2430.6667
$ws.Range("J102").Value = 2146.5
$ws.Range("K102").Value = 2430.6667
$ws.Range("L102").Value = 2146.5
$ws.Range("M102").Value = -808.6667000000002
$ws.Range("N102").Value = -5390.5

$ws.Range("H122").Value = 1826
$ws.Range("I122").Value = 1945.8572
$ws.Range("J122").Value = 1696.9231
$ws.Range("K122").Value = 5837.571599999999
$ws.Range("L122").Value = 5090.7693
$ws.Range("M122").Value = -3387.571599999999
$ws.Range("N122").Value = -9990.7693

$ws.Range("H126").Value = 20837148
$ws.Range("I126").Value = 3918.6667
$ws.Range("J126").Value = 83336830
$ws.Range("K126").Value = 11756.0001
$ws.Range("L126").Value = 250010490
$ws.Range("M126").Value = -9286.000100000001
$ws.Range("N126").Value = -250015430

$ws.Range("H132").Value = 29218.486
$ws.Range("I132").Value = 32745
$ws.Range("J132").Value = 5238.2
$ws.Range("K132").Value = 98235
$ws.Range("L132").Value = 15714.6
$ws.Range("M132").Value = -95705
$ws.Range("N132").Value = -20774.6

$ws.Range("H140").Value = 29840
$ws.Range("J140").Value = 29840
$ws.Range("L140").Value = 29840
$ws.Range("N140").Value = -40200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2303.1904
$ws.Range("I40").Value = 1882.5
$ws.Range("J40").Value = 2864.111
$ws.Range("K40").Value = 1882.5
$ws.Range("L40").Value = 2864.111
$ws.Range("M40").Value = -1746.5
$ws.Range("N40").Value = -3136.111

$ws.Range("H93").Value = 3212.389
$ws.Range("I93").Value = 3300.4285
$ws.Range("J93").Value = 3156.3635
$ws.Range("K93").Value = 3300.4285
$ws.Range("L93").Value = 3156.3635
$ws.Range("M93").Value = -2052.4285
$ws.Range("N93").Value = -5652.363499999999

$ws.Range("H100").Value = 3204.5789
$ws.Range("I100").Value = 2489
$ws.Range("J100").Value = 3725
$ws.Range("K100").Value = 2489
$ws.Range("L100").Value = 3725
$ws.Range("M100").Value = -1948
$ws.Range("N100").Value = -4807

$ws.Range("H132").Value = 6259.2666
$ws.Range("I132").Value = 7559.5
$ws.Range("J132").Value = 3658.8
$ws.Range("K132").Value = 22678.5
$ws.Range("L132").Value = 10976.4
$ws.Range("M132").Value = -20148.5
$ws.Range("N132").Value = -16036.4

$ws.Range("H139").Value = 79600
$ws.Range("J139").Value = 79600
$ws.Range("L139").Value = 79600
$ws.Range("N139").Value = -89880

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1499.4117
$ws.Range("I136").Value = 1023.75
$ws.Range("J136").Value = 2641
$ws.Range("K136").Value = 3071.25
$ws.Range("L136").Value = 7923
$ws.Range("M136").Value = -521.25
$ws.Range("N136").Value = -13023
